# "Generate Report for handback"
#
# The 3120aef5-5742-44c4-bfc8-f48d3381e7be source file has been handed
# back (it was previously "Ready for handoff"). Regenerate the
# Overview / zh-cn / de-de report rows so that:
#   - the 3120aef5... row reports "Handed back: in sync with en-US"
#     with its new handback file/datetime, and sorts above the
#     a1eefbe0... row (which was already handed back).
#   - the a1eefbe0... row keeps its existing (already handed-back)
#     data, just shifted down one row.
#   - the .localization-config row (row 4) is untouched.

$wb = $excel.ActiveWorkbook

function Set-CellText {
    param($ws, [string]$addr, [string]$text)
    $ws.Range($addr).Value = $text
}

function Set-LinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellText $wsOverview "A2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-CellText $wsOverview "B2" "Handed back: in sync with en-US"
Set-CellText $wsOverview "C2" "Handed back: in sync with en-US"

Set-CellText $wsOverview "A3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-CellText $wsOverview "B3" "Handed back: in sync with en-US"
Set-CellText $wsOverview "C3" "Handed back: in sync with en-US"

Set-LinkDisplay $wsOverview "`$A`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-LinkDisplay $wsOverview "`$A`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellText $wsZhCn "A2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-CellText $wsZhCn "B2" "Handed back: in sync with en-US"
Set-CellText $wsZhCn "C2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"
Set-CellText $wsZhCn "D2" "2016-01-17 03:20:01"
Set-CellText $wsZhCn "E2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-CellText $wsZhCn "F2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"
Set-CellText $wsZhCn "G2" "2016-01-17 03:20:44"
Set-CellText $wsZhCn "H2" "Include"

Set-CellText $wsZhCn "A3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-CellText $wsZhCn "B3" "Handed back: in sync with en-US"
Set-CellText $wsZhCn "C3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"
Set-CellText $wsZhCn "D3" "2016-01-17 03:18:19"
Set-CellText $wsZhCn "E3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-CellText $wsZhCn "F3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"
Set-CellText $wsZhCn "G3" "2016-01-17 03:19:00"
Set-CellText $wsZhCn "H3" "Include"

Set-LinkDisplay $wsZhCn "`$A`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-LinkDisplay $wsZhCn "`$C`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"
Set-LinkDisplay $wsZhCn "`$E`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-LinkDisplay $wsZhCn "`$F`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"

Set-LinkDisplay $wsZhCn "`$A`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-LinkDisplay $wsZhCn "`$C`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"
Set-LinkDisplay $wsZhCn "`$E`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-LinkDisplay $wsZhCn "`$F`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellText $wsDeDe "A2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-CellText $wsDeDe "B2" "Handed back: in sync with en-US"
Set-CellText $wsDeDe "C2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"
Set-CellText $wsDeDe "D2" "2016-01-17 03:20:11"
Set-CellText $wsDeDe "E2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-CellText $wsDeDe "F2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"
Set-CellText $wsDeDe "G2" "2016-01-17 03:21:00"
Set-CellText $wsDeDe "H2" "Include"

Set-CellText $wsDeDe "A3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-CellText $wsDeDe "B3" "Handed back: in sync with en-US"
Set-CellText $wsDeDe "C3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
Set-CellText $wsDeDe "D3" "2016-01-17 03:18:30"
Set-CellText $wsDeDe "E3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-CellText $wsDeDe "F3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
Set-CellText $wsDeDe "G3" "2016-01-17 03:19:17"
Set-CellText $wsDeDe "H3" "Include"

Set-LinkDisplay $wsDeDe "`$A`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-LinkDisplay $wsDeDe "`$C`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"
Set-LinkDisplay $wsDeDe "`$E`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
Set-LinkDisplay $wsDeDe "`$F`$2" "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"

Set-LinkDisplay $wsDeDe "`$A`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-LinkDisplay $wsDeDe "`$C`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
Set-LinkDisplay $wsDeDe "`$E`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
Set-LinkDisplay $wsDeDe "`$F`$3" "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
